# Update the "Förändrad" date column (C) for rows 2-32 from serial date
# 45175 (2023-09-06) to 45183 (2023-09-14), preserving existing formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C32").Value = 45183
